$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the merged-AHB column headers (row 1) ---
# "_old" columns now refer to the FV2410 formats-version, and "_new"
# columns now refer to the FV2504 formats-version.
$headerRange = $ws.Range("A1:U1")
[void]$headerRange.Replace("_old", "_FV2410")
[void]$headerRange.Replace("_new", "_FV2504")

# --- Freeze the header row ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn the data range into a proper table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"

[void]$ws.Range("A1").Select()
